$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.850.31"
$ws.Range("E2").Value = "  +1.02%  "

# Row 3
$ws.Range("D3").Value = "2.633.23"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.65"
$ws.Range("E5").Value = "  +4.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.05"
$ws.Range("E6").Value = "  +1.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.53%  "

# Row 10
$ws.Range("E10").Value = "  +5.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.350"
$ws.Range("E11").Value = "  +1.45%  "

# Row 12
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").Value = "3.094.83"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("D14").Value = "60.877.48"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15
$ws.Range("E15").Value = "  +2.26%  "

# Row 16
$ws.Range("E16").Value = "  +3.08%  "

# Row 17
$ws.Range("D17").Value = "2.637.64"
$ws.Range("E17").Value = "  +1.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.76"
$ws.Range("E18").Value = "  +0.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "353.25"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.62"
$ws.Range("E20").Value = "  +1.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  +2.05%  "

# Row 22
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.61"
$ws.Range("E23").Value = "  +2.15%  "

# Row 24
$ws.Range("E24").Value = "  +2.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  +1.43%  "

# Row 26
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0864"
$ws.Range("E27").Value = "  +3.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("E28").Value = "  +1.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.14"
$ws.Range("E30").Value = "  +7.25%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.48"
$ws.Range("E31").Value = "  +0.62%  "

# Row 32
$ws.Range("E32").Value = "  +4.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.41"
$ws.Range("E33").Value = "  -0.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  +4.62%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  +1.90%  "

# Row 36
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.925"
$ws.Range("E36").Value = "  +10.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.886"
$ws.Range("E37").Value = "  +2.41%  "

# Row 38
$ws.Range("E38").Value = "  +1.48%  "

# Row 39
$ws.Range("E39").Value = "  +2.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "306.46"
$ws.Range("E40").Value = "  +4.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.640"
$ws.Range("E41").Value = "  +3.99%  "

# Row 42
$ws.Range("E42").Value = "  +1.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0561"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.72"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.93"
$ws.Range("E46").Value = "  +3.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0238"
$ws.Range("E47").Value = "  +2.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.27"
$ws.Range("E48").Value = "  +8.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.35"
$ws.Range("E49").Value = "  +0.33%  "

# Row 50
$ws.Range("D50").Value = "1.979.35"
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.82"
$ws.Range("E51").Value = "  +2.85%  "
